$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

# New first search result (replaces the old "Images for yApATb" entry) -
# set the text first so the new shared strings are appended in the same
# order Excel would naturally allocate them.
$ws.Range("E3").Value = "Adidas One Grey W Gum4 Campus Core Black Footshop XgrXIq"
$ws.Range("F3").Value = "http://helper.extrapulpe.com/hcap-3-spanish.mdoc"

# Renumber the TestCaseID column and mark the existing cases as not-to-execute
$ws.Range("A3").Value = "No"
$ws.Range("B3").Value = "1"

$ws.Range("A4").Value = "No"
$ws.Range("B4").Value = "2"

$ws.Range("A5").Value = "No"
$ws.Range("B5").Value = "3"

# Add a new row 6 for the "windows" automation test case (BringToFront),
# copying the formatting from the row above it first.
$ws.Range("A5:G5").Copy($ws.Range("A6:G6"))

$ws.Range("A6").Value = "Yes"
$ws.Range("B6").Value = "4"
$ws.Range("C6").Value = "windows"
$ws.Range("D6").Value = "pass"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""

$ws.Range("A6").Select()
